# Add a new "Completed" book entry: "Beyond These Walls" by Tony Platt.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row right after the existing data (row 44 -> new row 45).
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "Beyond These Walls"
$ws.Cells.Item($newRow, 2).Value = "Tony Platt"
$ws.Cells.Item($newRow, 3).Value = 43913
$ws.Cells.Item($newRow, 4).Value = 43918
$ws.Cells.Item($newRow, 5).Value = "crime;prison reform;criminal justice;policing"
$ws.Cells.Item($newRow, 6).Value = "Hard Copy"
$ws.Cells.Item($newRow, 7).Value = "255 Pages"

# Match the existing date-format styling on the Start/Finish Date columns
# (copy formats from the previous row so the same cell style is reused).
$ws.Range("C" + $lastRow).Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)
$ws.Range("D" + $lastRow).Copy()
$ws.Range("D" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Mirror the author's selection after appending the row.
$ws.Range("A" + ($newRow + 1)).Select()
